# Update countries & provincias Spain
# - Swap order of "Ucrania" / "Indonesia" rows (Ucrania now sorts above Indonesia)
#   and refresh Ucrania's stats while Indonesia keeps its previous stats.
# - Update a handful of country/provincia stat rows with newer figures.
# - Bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 becomes Ucrania (fresh numbers), row 40 becomes Indonesia
# (keeps the numbers Indonesia previously had in row 39).
$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 9866
$ws.Range("C39").Value = 456
$ws.Range("D39").Value = 1103
$ws.Range("E39").Value = 8513
$ws.Range("F39").Value = 129
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 250

$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 9511
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 1254
$ws.Range("E40").Value = 7484
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 773

# Row 45 (Chequia): refreshed active/recovered counts
$ws.Range("D45").Value = 2960
$ws.Range("E45").Value = 4317

# Row 61 (Kazajistan): refreshed counts
$ws.Range("B61").Value = 3079
$ws.Range("C61").Value = 52
$ws.Range("E61").Value = 2280

# Row 98 (Kirguistan): refreshed deaths-today count
$ws.Range("F98").Value = 10

# Row 109 (Georgia): refreshed counts
$ws.Range("B109").Value = 517
$ws.Range("C109").Value = 6
$ws.Range("E109").Value = 343

# Update timestamp text
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 08:52"
